# Updated cryptos list values (Price and Volume(1h)) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): values are stored as literal text (e.g. "26.824.33"),
# so force a text number format before assignment to stop Excel from
# auto-converting numeric-looking strings into real numbers, then restore
# the default "Normal" style so no stray formatting is left behind.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.824.33'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.873.83'
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '301.45'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5371'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.8897'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.08171'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.875.49'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '93.42'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.319'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '14.85'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008535'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '26.863.99'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.990'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.289'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '146.53'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.733'
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '114.09'
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '4.721'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.619'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.09136'
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.8128'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.05015'
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.176'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.953'
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6025'
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.224'
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.626'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01958'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.070'
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.629'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '8.888'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '115.23'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.5118'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.1498'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '9.952'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '37.65'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06084'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '62.29'
$c.Style = "Normal"

# Volume(1h) column (E): percentage strings already contain "%" and
# padding spaces, so Excel keeps them as plain text automatically.
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("E13").Value = '  +2.80%  '
$ws.Range("E14").Value = '  -2.90%  '
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("E18").Value = '  -1.71%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("E20").Value = '  -1.54%  '
$ws.Range("E21").Value = '  -2.59%  '
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("E29").Value = '  -2.82%  '
$ws.Range("E30").Value = '  -5.43%  '
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("E34").Value = '  -4.90%  '
$ws.Range("E35").Value = '  -1.44%  '
$ws.Range("E36").Value = '  +4.59%  '
$ws.Range("E37").Value = '  -4.20%  '
$ws.Range("E38").Value = '  -4.76%  '
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("E49").Value = '  -2.73%  '
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("E51").Value = '  -2.98%  '
